$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
Write-Host "NotesMaster Shapes:" $nm.Shapes.Count
$theme = $nm.Theme
Write-Host "NotesMaster Theme:" $theme
Write-Host "ThemeVariants count:" $theme.ThemeVariants.Count
